## ISYS3001.docx - "Add files via upload"
##
## The paragraph that used to contain only the ellipsis character ("...")
## gets extra text appended to it:  " This is austin"
##
## In the authored version this text was typed with an IME active for part
## of the word "This", which is why Word split it into several runs (some
## of them carrying East-Asian font/language formatting). We reproduce the
## same run layout here: a run for the leading space, a run for "T", a run
## for "h", a run for "is" and a final run for " is austin".

$d = $word.ActiveDocument

# Find the paragraph whose whole text is the single ellipsis character.
$hit = $d.Content
$found = $hit.Find.Execute("…", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the paragraph containing the ellipsis character"
}

# $hit now covers exactly the "…" text; anchor new runs right after it
# (i.e. before the paragraph mark), so everything lands in that paragraph.
$pos = $hit.End

# --- " " --------------------------------------------------------------
$r = $d.Range($pos, $pos)
$r.InsertAfter(" ")
$pos = $pos + 1

# --- "T" (typed while an East-Asian input method / font was active) ---
$r = $d.Range($pos, $pos)
$r.InsertAfter("T")
$rT = $d.Range($pos, $pos + 1)
$rT.Font.Name = "SimSun"
$pos = $pos + 1

# --- "h" (typed back with the regular font) ----------------------------
$r = $d.Range($pos, $pos)
$r.InsertAfter("h")
$pos = $pos + 1

# --- "is" (East-Asian font again) --------------------------------------
$r = $d.Range($pos, $pos)
$r.InsertAfter("is")
$rIs = $d.Range($pos, $pos + 2)
$rIs.Font.Name = "SimSun"
$pos = $pos + 2

# --- " is austin" (rest typed normally) ---------------------------------
$r = $d.Range($pos, $pos)
$r.InsertAfter(" is austin")
$pos = $pos + 10
